$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the two descriptive labels in column A (shared-string text edits)
$ws.Range("A6").Value = "Willekeurige afschrijving investeringen"
$ws.Range("A10").Value = "Vervallen vrijstelling MRB auto's>25jr"

# Move the active selection from B1 to A10
$ws.Range("A10").Select()
